# Generate Report for Handoff
# Refresh the "Latest Handoff Datetime" / "Latest HO Xliff Generate Date" values
# for the f9fa95ac-a9e9-41dc-9a83-0c55b48da1ef.md file after a fresh handoff run.

$wb = $excel.ActiveWorkbook

# de-de locale sheet: the handoff xliff for f9fa95ac was (re)generated at 06:12:01
$deDeSheet = $wb.Worksheets.Item("de-de")
$deDeSheet.Range("H5").Value = "2016-10-27 06:12:01"

# zh-cn locale sheet: the handoff xliff for f9fa95ac was (re)generated at 06:11:50
$zhCnSheet = $wb.Worksheets.Item("zh-cn")
$zhCnSheet.Range("H5").Value = "2016-10-27 06:11:50"

# Overview sheet: roll up the latest of the two per-locale handoff datetimes
$overviewSheet = $wb.Worksheets.Item("Overview")
$overviewSheet.Range("G5").Value = "2016-10-27 06:12:01"
